$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dbo_counmet")

# Insert a new row above current row 2, shifting existing rows (2..39) down to (3..40)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the "Not applicable" lookup entry
$ws.Cells.Item(2, 1).Value = -1
$ws.Cells.Item(2, 2).Value = "Not applicable"

# Update the sheet selection to match the authored state
$ws.Range("A3").Select()

# Keep the workbook-level defined name in sync with the new data extent
$wb.Names.Item("dbo_counmet").RefersTo = "=dbo_counmet!`$A`$1:`$C`$40"
